$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column C slightly to fit the new "Linked List" pattern text
$ws.Columns.Item(3).ColumnWidth = 15

# --- Row 10: 191. Number of 1 Bits (Easy / Bit Manipulation) ---
$ws.Range("A10").Value = "191. Number of 1 Bits"
$ws.Range("B10").Value = "Easy"
$ws.Range("B10").Interior.Color = 5287936
$ws.Range("C10").Value = "Bit Manipulation"
$ws.Range("D10").Value = "Bit Shifting: While n != 0, count += (n&1), unsigned right shift (>>>) by 1. Optimal solution is n = n & (n-1) without shifting."
$ws.Hyperlinks.Add($ws.Range("E10"), "https://leetcode.com/problems/number-of-1-bits/solutions/55099/simple-java-solution-bit-shifting/ ")
$ws.Range("E10").Style = "Hyperlink"

# --- Row 11: 19. Remove Nth Node From End of List (Medium / Linked List) ---
$ws.Range("A11").Value = "19. Remove Nth Node From End of List"
$ws.Range("B11").Value = "Medium"
$ws.Range("B11").Interior.Color = 49407
$ws.Range("C11").Value = "Linked List"
$ws.Range("D11").Value = "Fast and Slow pointers, move fast n ahead, slow.next = slow.next.next"
$ws.Hyperlinks.Add($ws.Range("E11"), "https://leetcode.com/problems/remove-nth-node-from-end-of-list/solutions/1164542/js-python-java-c-easy-two-pointer-solution-w-explanation/ ")
$ws.Range("E11").Style = "Hyperlink"

# --- Row 12: 90. Subsets 2 (Medium / Backtracking) ---
$ws.Range("A12").Value = "90. Subsets 2"
$ws.Range("B12").Value = "Medium"
$ws.Range("B12").Interior.Color = 49407
$ws.Range("C12").Value = "Backtracking"
$ws.Hyperlinks.Add($ws.Range("E12"), "https://leetcode.com/problems/subsets-ii/solutions/388566/subsets-i-ii-java-solution-with-detailed-explanation-and-comments-recursion-iteration/ ")
$ws.Range("E12").Style = "Hyperlink"
$ws.Range("D12").Value = "Power Set, skip iterations with duplicate previous element."

# Extend the table to cover the new rows
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:E12"))

$ws.Range("D13").Select()
